# Rename several "pf_*" result-column headers across all worksheets of the
# workbook. Sheets 1-8 ("LLL_*") use a 17-column layout (A:Q); sheets 9-32
# ("LL_*", "LLG_*", "LG_*") use a 43-column layout (A:AQ) with per-phase
# (a/b/c) columns. Only header row 1 is affected.

$wb = $excel.ActiveWorkbook

# Mapping used for the narrow (A:Q) header layout - sheets 1-8.
$narrowMap = @{
    "L1" = "pf_ikss_from_degree"
    "M1" = "pf_ikss_to_degree"
    "P1" = "pf_va_from_degree"
    "Q1" = "pf_va_to_degree"
}

# Mapping used for the wide (A:AQ) header layout - sheets 9-32.
$wideMap = @{
    "T1"  = "pf_q_a_from_mvar"
    "U1"  = "pf_q_b_from_mvar"
    "V1"  = "pf_q_c_from_mvar"
    "W1"  = "pf_q_a_to_mvar"
    "X1"  = "pf_q_b_to_mvar"
    "Y1"  = "pf_q_c_to_mvar"
    "Z1"  = "pf_ikss_a_from_degree"
    "AA1" = "pf_ikss_b_from_degree"
    "AB1" = "pf_ikss_c_from_degree"
    "AC1" = "pf_ikss_a_to_degree"
    "AD1" = "pf_ikss_b_to_degree"
    "AE1" = "pf_ikss_c_to_degree"
    "AG1" = "pf_vm_b_from_pu"
    "AH1" = "pf_vm_c_from_pu"
    "AI1" = "pf_vm_a_to_pu"
    "AJ1" = "pf_vm_b_to_pu"
    "AK1" = "pf_vm_c_to_pu"
    "AL1" = "pf_va_a_from_degree"
    "AM1" = "pf_va_b_from_degree"
    "AN1" = "pf_va_c_from_degree"
    "AO1" = "pf_va_a_to_degree"
    "AP1" = "pf_va_b_to_degree"
    "AQ1" = "pf_va_c_to_degree"
}

foreach ($ws in $wb.Worksheets) {
    # Distinguish the two header layouts by the sheet's used-column count
    # (17 columns => narrow "LLL_*" sheets, 43 columns => wide sheets with
    # per-phase a/b/c columns) rather than by sheet name, since that is what
    # actually determines which header cells exist.
    $colCount = $ws.UsedRange.Columns.Count

    if ($colCount -gt 17) {
        $map = $wideMap
    } else {
        $map = $narrowMap
    }

    foreach ($addr in $map.Keys) {
        $ws.Range($addr).Value = $map[$addr]
    }
}
